$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1930666666666666
$ws.Range("H2").Value = 0.5791999999999999
$ws.Range("I2").Value = 0.01292026122037801
$ws.Range("J2").Value = 0.01292026122037801
$ws.Range("M2").Value = 0.4652636666666667
$ws.Range("N2").Value = 1.395791
$ws.Range("O2").Value = 0.02604271297411062
$ws.Range("P2").Value = 0.02604271297411062
$ws.Range("Q2").Value = 0.08982690524444444
$ws.Range("R2").Value = 0.8084421472
$ws.Range("S2").Value = 0.0003364786545128367
$ws.Range("T2").Value = 0.0003364786545128366
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1930666666666666
$ws.Range("H3").Value = 0.5791999999999999
$ws.Range("I3").Value = 0.01292026122037801
$ws.Range("J3").Value = 0.01292026122037801
$ws.Range("O3").Value = 0.09971126509087273
$ws.Range("P3").Value = 0.09971126509087272
$ws.Range("Q3").Value = 0.3439255491555555
$ws.Range("R3").Value = 3.095329942399999
$ws.Range("S3").Value = 0.001288295591588434
$ws.Range("T3").Value = 0.001288295591588434
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1930666666666666
$ws.Range("H4").Value = 0.5791999999999999
$ws.Range("I4").Value = 0.01292026122037801
$ws.Range("J4").Value = 0.01292026122037801
$ws.Range("M4").Value = 15.618761
$ws.Range("N4").Value = 46.856283
$ws.Range("O4").Value = 0.8742460219350168
$ws.Range("P4").Value = 0.8742460219350167
$ws.Range("Q4").Value = 3.015462123733333
$ws.Range("R4").Value = 27.13915911359999
$ws.Range("S4").Value = 0.01129548697427674
$ws.Range("T4").Value = 0.01129548697427674
$ws.Range("I5").Value = 0.04457951877603724
$ws.Range("J5").Value = 0.04457951877603725
$ws.Range("M5").Value = 0.4652636666666667
$ws.Range("N5").Value = 1.395791
$ws.Range("O5").Value = 0.02604271297411062
$ws.Range("P5").Value = 0.02604271297411062
$ws.Range("Q5").Value = 0.3099349262863333
$ws.Range("R5").Value = 2.789414336577
$ws.Range("S5").Value = 0.001160971612008313
$ws.Range("T5").Value = 0.001160971612008313
$ws.Range("I6").Value = 0.04457951877603724
$ws.Range("J6").Value = 0.04457951877603725
$ws.Range("O6").Value = 0.09971126509087273
$ws.Range("P6").Value = 0.09971126509087272
$ws.Range("S6").Value = 0.004445080214300987
$ws.Range("T6").Value = 0.004445080214300988
$ws.Range("I7").Value = 0.04457951877603724
$ws.Range("J7").Value = 0.04457951877603725
$ws.Range("M7").Value = 15.618761
$ws.Range("N7").Value = 46.856283
$ws.Range("O7").Value = 0.8742460219350168
$ws.Range("P7").Value = 0.8742460219350167
$ws.Range("Q7").Value = 10.404422021389
$ws.Range("R7").Value = 93.639798192501
$ws.Range("S7").Value = 0.03897346694972795
$ws.Range("T7").Value = 0.03897346694972795
$ws.Range("G8").Value = 14.08372266666667
$ws.Range("H8").Value = 42.251168
$ws.Range("I8").Value = 0.9425002200035847
$ws.Range("J8").Value = 0.9425002200035848
$ws.Range("M8").Value = 0.4652636666666667
$ws.Range("N8").Value = 1.395791
$ws.Range("O8").Value = 0.02604271297411062
$ws.Range("P8").Value = 0.02604271297411062
$ws.Range("Q8").Value = 6.552644448209778
$ws.Range("R8").Value = 58.973800033888
$ws.Range("S8").Value = 0.02454526270758947
$ws.Range("T8").Value = 0.02454526270758947
$ws.Range("G9").Value = 14.08372266666667
$ws.Range("H9").Value = 42.251168
$ws.Range("I9").Value = 0.9425002200035847
$ws.Range("J9").Value = 0.9425002200035848
$ws.Range("O9").Value = 0.09971126509087273
$ws.Range("P9").Value = 0.09971126509087272
$ws.Range("Q9").Value = 25.08849474596622
$ws.Range("R9").Value = 225.796452713696
$ws.Range("S9").Value = 0.09397788928498331
$ws.Range("T9").Value = 0.09397788928498331
$ws.Range("G10").Value = 14.08372266666667
$ws.Range("H10").Value = 42.251168
$ws.Range("I10").Value = 0.9425002200035847
$ws.Range("J10").Value = 0.9425002200035848
$ws.Range("M10").Value = 15.618761
$ws.Range("N10").Value = 46.856283
$ws.Range("O10").Value = 0.8742460219350168
$ws.Range("P10").Value = 0.8742460219350167
$ws.Range("Q10").Value = 219.9702983209493
$ws.Range("R10").Value = 1979.732684888544
$ws.Range("S10").Value = 0.823977068011012
$ws.Range("T10").Value = 0.823977068011012